$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19
$ws.Cells.Item($row, 1).Value = 11.81
$ws.Cells.Item($row, 2).Value = 11.21
$ws.Cells.Item($row, 3).Value = 1019
$ws.Cells.Item($row, 4).Value = 83
$ws.Cells.Item($row, 5).Value = "few clouds"
$ws.Cells.Item($row, 6).Value = 20
$ws.Cells.Item($row, 7).Value = "Lisbon"
$ws.Cells.Item($row, 8).Value = 18.504
$ws.Cells.Item($row, 9).Value = 20
$ws.Cells.Item($row, 10).Value = "19:55:41 02-12-2025"
